# Edit script: shift slide content (slide N gets the title/body that used to
# belong conceptually to slide N+1's topic), then drop the now-redundant
# trailing slide (old slide 7 "提交方式").
#
# Each content slide (2-7) has exactly two shapes:
#   Shapes.Item(1) -> Title placeholder (single paragraph)
#   Shapes.Item(2) -> Content placeholder; Paragraphs(1) is an empty leading
#                     paragraph, Paragraphs(2) and Paragraphs(3) hold the two
#                     bullet lines (slide 7 only has one bullet line).
#
# NOTE: setting Paragraphs(i).Text directly to a string that shares a common
# prefix with the existing text causes the runtime to keep the shared prefix
# as its own run and append a new run for the changed suffix (to preserve
# formatting). To always end up with a single clean <a:r><a:t> run that
# matches the original authoring, we first blank the paragraph text and then
# set the desired final text.

function Set-ParaText($para, $text) {
    $para.Text = ""
    $para.Text = $text
}

$p = $ppt.ActivePresentation

# --- Slide 2 ---
$s = $p.Slides.Item(2)
Set-ParaText $s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1) "使用 python-pptx 自动生成 PowerPoint 文件"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr.Paragraphs(2) "使用 python-pptx 包创建演示文稿"
Set-ParaText $tr.Paragraphs(3) "支持文本、图片、表格和图表的插入"

# --- Slide 3 ---
$s = $p.Slides.Item(3)
Set-ParaText $s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1) "Gradio 搭建 ChatBot"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr.Paragraphs(2) "构建图形化用户界面"
Set-ParaText $tr.Paragraphs(3) "将用户输入转化为 ChatPPT PowerPoint 输入格式（Markdown）"

# --- Slide 4 ---
$s = $p.Slides.Item(4)
Set-ParaText $s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1) "ChatBot System Prompt"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr.Paragraphs(2) "使用 ChatPPT v0.2 prompts/formatter.txt 文件"
Set-ParaText $tr.Paragraphs(3) "鼓励自行创作和优化"

# --- Slide 5 ---
$s = $p.Slides.Item(5)
Set-ParaText $s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1) "整合主流程（可选）"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr.Paragraphs(2) "支持聊天输入"
Set-ParaText $tr.Paragraphs(3) "自动生成 PowerPoint 文件作为输出"

# --- Slide 6 ---
$s = $p.Slides.Item(6)
Set-ParaText $s.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1) "作业提交方式"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr.Paragraphs(2) "修改代码文件链接复制粘贴至评论框"
Set-ParaText $tr.Paragraphs(3) "点击提交按钮完成作业"

# --- Remove old slide 7 ("提交方式") entirely ---
$p.Slides.Item(7).Delete()
